$d = $word.ActiveDocument

# --- Locate the three split points around "using a framework" / "in python" ---

# Point A: right after " a neural network" (i.e. right before " using a framework")
$rng = $d.Content
$rng.Find.Execute("a neural network", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pointA = $rng.End

# Point B: right after " using a framework" (i.e. right before " in python")
$rng = $d.Content
$rng.Find.Execute("using a framework", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pointB = $rng.End

# Point C: right after " in python" (i.e. right before the following space/"incorporating")
$rng = $d.Content
$rng.Find.Execute("in python", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pointC = $rng.End

# --- Fence each edit region with temporary bookmarks so that this runtime's ---
# --- run-coalescing (which merges same-formatted runs touched by an edit)  ---
# --- cannot spread past the boundaries we care about. Bookmarks re-anchor  ---
# --- themselves as the surrounding text is edited, so we read their live  ---
# --- Range back out before each subsequent edit instead of reusing stale  ---
# --- integer offsets.                                                      ---
$d.Bookmarks.Add("zzzTempA", $d.Range($pointA, $pointA))
$d.Bookmarks.Add("zzzTempB", $d.Range($pointB, $pointB))
$d.Bookmarks.Add("zzzTempC", $d.Range($pointC, $pointC))

# Edit 1: " using a framework" -> " "  (drop "using a framework", keep one space)
$rA = $d.Bookmarks.Item("zzzTempA").Range
$rB = $d.Bookmarks.Item("zzzTempB").Range
$d.Range($rA.Start, $rB.Start).Text = " "

# Edit 2: " in python" -> "in python" (drop the leading space)
$rB = $d.Bookmarks.Item("zzzTempB").Range
$rC = $d.Bookmarks.Item("zzzTempC").Range
$d.Range($rB.Start, $rC.Start).Text = "in python"

# --- Turn the (now zero-length) middle fence-post into the real "_GoBack" ---
# --- bookmark Word leaves at the location of the most recent edit, and    ---
# --- drop the two outer fence posts we only needed during the edit.      ---
$mid = $d.Bookmarks.Item("zzzTempB").Range
$d.Bookmarks.Item("zzzTempA").Delete()
$d.Bookmarks.Item("zzzTempC").Delete()

# Adding a bookmark named "_GoBack" automatically relocates/replaces any
# existing one elsewhere in the document (there is only ever one).
$d.Bookmarks.Add("_GoBack", $mid)
$d.Bookmarks.Item("zzzTempB").Delete()
